# Bitacora.docx – "correccion en la documentacion"
#
# 1) The 4th ("Flor Machado") and 5th ("Elías Peregrina") data columns of
#    the log table are resized: 2306 -> 2172 dxa and 1795 -> 1929 dxa
#    (115.3pt -> 108.6pt, 89.75pt -> 96.45pt). Setting Column.Width
#    updates both the <w:tblGrid> entry and every row's <w:tcW> for that
#    column in one shot.
# 2) In the last row of the table (row 8 = the most recent log entry),
#    both people's task cells get an extra "Documentación" bullet
#    appended after their existing "Escribir código" bullet, re-using
#    the same list style/numbering already on that paragraph.

$table = $word.ActiveDocument.Tables.Item(1)
$table.Columns.Item(4).Width = 108.6

$table = $word.ActiveDocument.Tables.Item(1)
$table.Columns.Item(5).Width = 96.45

# Column 4 ("Flor Machado") of the last row: add "Documentación" bullet.
$table = $word.ActiveDocument.Tables.Item(1)
$cell = $table.Cell(8, 4)
$cell.Range.Paragraphs.Item(1).Range.InsertAfter([char]13 + "Documentación")

# Column 5 ("Elías Peregrina") of the last row: add "Documentación" bullet.
$table = $word.ActiveDocument.Tables.Item(1)
$cell = $table.Cell(8, 5)
$cell.Range.Paragraphs.Item(1).Range.InsertAfter([char]13 + "Documentación")
